$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pop times")

# Replace "iterations" with "models" in the relevant cells
$ws.Range("E1").Value = "Number of different models per scenario"
$ws.Range("D7").Value = "Total models"

# Update the active selection to match the final cursor position
$ws.Range("D8").Select()
